# Fix gaussian and kmean: relabel the "Wood Class" (k-means cluster) column
# Column E holds the cluster id produced by the k-means step. The cluster
# labels got shuffled, so remap the old label -> new label consistently:
#   0 -> 2
#   1 -> 0
#   2 -> 5
#   3 -> 3 (unchanged)
#   4 -> 4 (unchanged)
#   5 -> 1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{ 0 = 2; 1 = 0; 2 = 5; 3 = 3; 4 = 4; 5 = 1 }

$lastRow = $ws.Cells.Item($ws.Rows.Count, 5).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $old = $cell.Value2
    if ($null -ne $old) {
        $oldInt = [int]$old
        if ($map.ContainsKey($oldInt)) {
            $cell.Value2 = $map[$oldInt]
        }
    }
}
